$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 119.0815153333333
    "H2" = 357.244546
    "I2" = 0.431812569872284
    "J2" = 0.4318125698722839
    "M2" = 45.90594266666667
    "N2" = 137.717828
    "O2" = 0.3954672001633582
    "P2" = 0.3954672001633583
    "Q2" = 5466.549215551788
    "R2" = 49198.94293996609
    "S2" = 0.1707677080027366
    "T2" = 0.1707677080027367
    "G3" = 119.0815153333333
    "H3" = 357.244546
    "I3" = 0.431812569872284
    "J3" = 0.4318125698722839
    "O3" = 0.3484294080560655
    "P3" = 0.3484294080560656
    "Q3" = 4816.345088789334
    "R3" = 43347.10579910401
    "S3" = 0.1504561981117684
    "T3" = 0.1504561981117684
    "G4" = 119.0815153333333
    "H4" = 357.244546
    "I4" = 0.431812569872284
    "J4" = 0.4318125698722839
    "M4" = 12.761795
    "N4" = 38.28538500000001
    "O4" = 0.1099393900775594
    "P4" = 0.1099393900775594
    "Q4" = 1519.693886973357
    "R4" = 13677.24498276021
    "S4" = 0.0474732105595824
    "T4" = 0.04747321055958241
    "G5" = 119.0815153333333
    "H5" = 357.244546
    "I5" = 0.431812569872284
    "J5" = 0.4318125698722839
    "M5" = 16.966758
    "N5" = 50.900274
    "O5" = 0.1461640017030168
    "P5" = 0.1461640017030168
    "Q5" = 2020.427252933956
    "R5" = 18183.8452764056
    "S5" = 0.06311545319819659
    "T5" = 0.06311545319819659
    "I6" = 0.4460879372303943
    "J6" = 0.4460879372303942
    "M6" = 45.90594266666667
    "N6" = 137.717828
    "O6" = 0.3954672001633582
    "P6" = 0.3954672001633583
    "Q6" = 5647.268823265553
    "R6" = 50825.41940938996
    "S6" = 0.1764131475631519
    "T6" = 0.1764131475631519
    "I7" = 0.4460879372303943
    "J7" = 0.4460879372303942
    "O7" = 0.3484294080560655
    "P7" = 0.3484294080560656
    "Q7" = 4975.569484425232
    "S7" = 0.1554301559101376
    "T7" = 0.1554301559101376
    "I8" = 0.4460879372303943
    "J8" = 0.4460879372303942
    "M8" = 12.761795
    "N8" = 38.28538500000001
    "O8" = 0.1099393900775594
    "P8" = 0.1099393900775594
    "Q8" = 1569.93371328234
    "R8" = 14129.40341954106
    "S8" = 0.04904263574006616
    "T8" = 0.04904263574006616
    "I9" = 0.4460879372303943
    "J9" = 0.4460879372303942
    "M9" = 16.966758
    "N9" = 50.900274
    "O9" = 0.1461640017030168
    "P9" = 0.1461640017030168
    "Q9" = 2087.220911267016
    "R9" = 18784.98820140314
    "S9" = 0.0652019980170386
    "T9" = 0.0652019980170386
    "G10" = 33.50679633333333
    "H10" = 100.520389
    "I10" = 0.1215021138451521
    "J10" = 0.121502113845152
    "M10" = 45.90594266666667
    "N10" = 137.717828
    "O10" = 0.3954672001633582
    "P10" = 0.3954672001633583
    "Q10" = 1538.161071421677
    "R10" = 13843.44964279509
    "S10" = 0.04805010077627189
    "T10" = 0.04805010077627189
    "G11" = 33.50679633333333
    "H11" = 100.520389
    "I11" = 0.1215021138451521
    "J11" = 0.121502113845152
    "O11" = 0.3484294080560655
    "P11" = 0.3484294080560656
    "Q11" = 1355.208602354264
    "R11" = 12196.87742118837
    "S11" = 0.04233490960462702
    "T11" = 0.04233490960462702
    "G12" = 33.50679633333333
    "H12" = 100.520389
    "I12" = 0.1215021138451521
    "J12" = 0.121502113845152
    "M12" = 12.761795
    "N12" = 38.28538500000001
    "O12" = 0.1099393900775594
    "P12" = 0.1099393900775594
    "Q12" = 427.6068659127517
    "R12" = 3848.461793214765
    "S12" = 0.0133578682892702
    "T12" = 0.01335786828927021
    "G13" = 33.50679633333333
    "H13" = 100.520389
    "I13" = 0.1215021138451521
    "J13" = 0.121502113845152
    "M13" = 16.966758
    "N13" = 50.900274
    "O13" = 0.1461640017030168
    "P13" = 0.1461640017030168
    "Q13" = 568.501704742954
    "R13" = 5116.515342686585
    "S13" = 0.01775923517498295
    "T13" = 0.01775923517498295
    "E14" = 2
    "F14" = 0.6666666666666666
    "G14" = 0.16474
    "H14" = 0.49422
    "I14" = 0.000597379052169715
    "J14" = 0.000597379052169715
    "M14" = 45.90594266666667
    "N14" = 137.717828
    "O14" = 0.3954672001633582
    "P14" = 0.3954672001633583
    "Q14" = 7.562544994906667
    "R14" = 68.06290495416
    "S14" = 0.0002362438211977979
    "T14" = 0.000236243821197798
    "E15" = 2
    "F15" = 0.6666666666666666
    "G15" = 0.16474
    "H15" = 0.49422
    "I15" = 0.000597379052169715
    "J15" = 0.000597379052169715
    "O15" = 0.3484294080560655
    "P15" = 0.3484294080560656
    "Q15" = 6.663038236506667
    "R15" = 59.96734412855999
    "S15" = 0.0002081444295325873
    "T15" = 0.0002081444295325873
    "E16" = 2
    "F16" = 0.6666666666666666
    "G16" = 0.16474
    "H16" = 0.49422
    "I16" = 0.000597379052169715
    "J16" = 0.000597379052169715
    "M16" = 12.761795
    "N16" = 38.28538500000001
    "O16" = 0.1099393900775594
    "P16" = 0.1099393900775594
    "Q16" = 2.1023781083
    "R16" = 18.9214029747
    "S16" = 0.000065675488640649
    "T16" = 0.00006567548864064903
    "E17" = 2
    "F17" = 0.6666666666666666
    "G17" = 0.16474
    "H17" = 0.49422
    "I17" = 0.000597379052169715
    "J17" = 0.000597379052169715
    "M17" = 16.966758
    "N17" = 50.900274
    "O17" = 0.1461640017030168
    "P17" = 0.1461640017030168
    "Q17" = 2.79510371292
    "R17" = 25.15593341628
    "S17" = 0.0000873153127986808
    "T17" = 0.00008731531279868081
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
